$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1894409937888199
$ws.Range("C2").Value = 0.562111801242236
$ws.Range("J2").Value = 0.03105590062111801
$ws.Range("P2").Value = 0.1118012422360248
$ws.Range("S2").Value = 0.1055900621118012
$ws.Range("B3").Value = 0.01030927835051546
$ws.Range("C3").Value = 0.06185567010309279
$ws.Range("J3").Value = 0.06701030927835051
$ws.Range("P3").Value = 0.6958762886597938
$ws.Range("S3").Value = 0.1649484536082474
$ws.Range("J4").Value = 0.07142857142857142
$ws.Range("P4").Value = 0.5714285714285714
$ws.Range("S4").Value = 0.3571428571428572
$ws.Range("B6").Value = 0.06584362139917696
$ws.Range("D6").Value = 0.02469135802469136
$ws.Range("E6").Value = 0.00411522633744856
$ws.Range("F6").Value = 0.06995884773662552
$ws.Range("J6").Value = 0.2674897119341564
$ws.Range("O6").Value = 0.01234567901234568
$ws.Range("Q6").Value = 0.1604938271604938
$ws.Range("R6").Value = 0.05761316872427984
$ws.Range("S6").Value = 0.3374485596707819
$ws.Range("B7").Value = 0.095
$ws.Range("D7").Value = 0.005
$ws.Range("F7").Value = 0.04
$ws.Range("J7").Value = 0.08
$ws.Range("O7").Value = 0.03
$ws.Range("Q7").Value = 0.24
$ws.Range("R7").Value = 0.1
$ws.Range("S7").Value = 0.41
$ws.Range("B8").Value = 0.1091269841269841
$ws.Range("D8").Value = 0.01587301587301587
$ws.Range("F8").Value = 0.07142857142857142
$ws.Range("J8").Value = 0.08134920634920635
$ws.Range("O8").Value = 0.01587301587301587
$ws.Range("Q8").Value = 0.2043650793650794
$ws.Range("R8").Value = 0.0992063492063492
$ws.Range("S8").Value = 0.4027777777777778
$ws.Range("B9").Value = 0.1191709844559585
$ws.Range("D9").Value = 0.02590673575129534
$ws.Range("F9").Value = 0.07253886010362694
$ws.Range("J9").Value = 0.08290155440414508
$ws.Range("O9").Value = 0.01036269430051814
$ws.Range("Q9").Value = 0.155440414507772
$ws.Range("R9").Value = 0.1243523316062176
$ws.Range("S9").Value = 0.4093264248704663
$ws.Range("B10").Value = 0.09909281228192603
$ws.Range("D10").Value = 0.01744591765526867
$ws.Range("E10").Value = 0.001395673412421493
$ws.Range("F10").Value = 0.05373342637822749
$ws.Range("J10").Value = 0.09281228192602931
$ws.Range("O10").Value = 0.01186322400558269
$ws.Range("Q10").Value = 0.2226099092812282
$ws.Range("R10").Value = 0.1165387299371947
$ws.Range("S10").Value = 0.3845080251221214
$ws.Range("F11").Value = 0.00308641975308642
$ws.Range("G11").Value = 0.1327160493827161
$ws.Range("J11").Value = 0.09876543209876543
$ws.Range("K11").Value = 0.2067901234567901
$ws.Range("L11").Value = 0.5339506172839507
$ws.Range("S11").Value = 0.02469135802469136
$ws.Range("G12").Value = 0.7288135593220338
$ws.Range("J12").Value = 0.2033898305084746
$ws.Range("K12").Value = 0.01694915254237288
$ws.Range("L12").Value = 0.03389830508474576
$ws.Range("S12").Value = 0.01694915254237288
$ws.Range("G13").Value = 0.6792452830188679
$ws.Range("J13").Value = 0.3018867924528302
$ws.Range("S13").Value = 0.01886792452830189
$ws.Range("F15").Value = 0.03018867924528302
$ws.Range("H15").Value = 0.169811320754717
$ws.Range("I15").Value = 0.04905660377358491
$ws.Range("J15").Value = 0.3773584905660378
$ws.Range("K15").Value = 0.04528301886792453
$ws.Range("M15").Value = 0.01132075471698113
$ws.Range("N15").Value = 0.003773584905660377
$ws.Range("O15").Value = 0.09056603773584905
$ws.Range("S15").Value = 0.2226415094339623
$ws.Range("F16").Value = 0.03191489361702127
$ws.Range("H16").Value = 0.1702127659574468
$ws.Range("I16").Value = 0.07446808510638298
$ws.Range("J16").Value = 0.3563829787234042
$ws.Range("K16").Value = 0.1170212765957447
$ws.Range("M16").Value = 0.02659574468085106
$ws.Range("O16").Value = 0.1063829787234043
$ws.Range("S16").Value = 0.1170212765957447
$ws.Range("F17").Value = 0.03171641791044776
$ws.Range("H17").Value = 0.1567164179104478
$ws.Range("I17").Value = 0.07649253731343283
$ws.Range("J17").Value = 0.4533582089552239
$ws.Range("K17").Value = 0.08768656716417911
$ws.Range("M17").Value = 0.02425373134328358
$ws.Range("O17").Value = 0.06716417910447761
$ws.Range("S17").Value = 0.1026119402985075
$ws.Range("F18").Value = 0.04379562043795621
$ws.Range("H18").Value = 0.1642335766423358
$ws.Range("I18").Value = 0.072992700729927
$ws.Range("J18").Value = 0.4343065693430657
$ws.Range("K18").Value = 0.0948905109489051
$ws.Range("M18").Value = 0.0218978102189781
$ws.Range("O18").Value = 0.06934306569343066
$ws.Range("S18").Value = 0.09854014598540146
$ws.Range("F19").Value = 0.01217765042979943
$ws.Range("H19").Value = 0.2134670487106017
$ws.Range("I19").Value = 0.07664756446991404
$ws.Range("J19").Value = 0.3861031518624642
$ws.Range("K19").Value = 0.1017191977077364
$ws.Range("M19").Value = 0.02077363896848138
$ws.Range("N19").Value = 0.001432664756446991
$ws.Range("O19").Value = 0.07163323782234957
$ws.Range("S19").Value = 0.1160458452722063
